# Add data for 2021-12-08
# - Rename sheet tab from "Through 2021-11-29" to "Through 2021-11-30"
# - Update the "November (through 11-29/30)" row (row 13) and the
#   "Total" row (row 14) with the latest counts / rates for every year
#   column (2015-2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-11-30"

# Update the row label.
$ws.Range("A13").Value = "November (through 11-30)"

# --- Row 13 ("November (through ...)") ---------------------------------
# 2015
$ws.Range("C13").Value = 32
$ws.Range("D13").Value = 0.0303
# 2016
$ws.Range("F13").Value = 69
$ws.Range("G13").Value = 0.1039
# 2017
$ws.Range("I13").Value = 109
$ws.Range("J13").Value = 0.018
# 2019
$ws.Range("O13").Value = 46
$ws.Range("P13").Value = 0.1154
# 2020
$ws.Range("R13").Value = 197
$ws.Range("S13").Value = 0.0483
# 2021
$ws.Range("U13").Value = 198
$ws.Range("V13").Value = 0.0198

# --- Row 14 ("Total") ---------------------------------------------------
# 2015
$ws.Range("C14").Value = 258
$ws.Range("D14").Value = 0.1134
# 2016
$ws.Range("F14").Value = 503
$ws.Range("G14").Value = 0.1066
# 2017
$ws.Range("I14").Value = 758
$ws.Range("J14").Value = 0.0767
# 2019
$ws.Range("O14").Value = 480
$ws.Range("P14").Value = 0.1011
# 2020
$ws.Range("R14").Value = 1200
$ws.Range("S14").Value = 0.0506
# 2021
$ws.Range("U14").Value = 1549
$ws.Range("V14").Value = 0.0584
